# examen git numero 2
# Adds a new set of Git quiz question/answer rows to the "GIT" sheet,
# continuing the existing alternating "question" (Good/green) /
# "answer" (bold) row pattern, and updates the sheet's view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GIT")

# The sheet already uses two visual row styles:
#   - A1-style (green "Good"/"Buena") for the question/prompt rows
#   - A2-style (bold dark-gray)        for the answer rows
# Re-use those exact formats (copy/paste-special formats) instead of
# assigning named styles, so no duplicate style/font entries get created.

function Add-QA($Row, $Text, $Kind) {
    # $Kind: "Q" = question/prompt row (green "Good" style, like A1)
    #        "A" = answer row (bold style, like A2)
    $cell = $ws.Range("A$Row")
    $cell.Value = $Text

    if ($Kind -eq "Q") {
        $ws.Range("A1").Copy()
    } else {
        $ws.Range("A2").Copy()
    }
    $cell.PasteSpecial(-4122)
}

Add-QA 13 "Si al momento de querer aplicar un stash el working directory tiene cambios en los mismos archivos que dicho stash modifica:" "Q"
Add-QA 14 " Se aplicará pero los archivos quedarán en estado de conflicto" "A"
Add-QA 15 "¿En qué se parecen un tag y un branch?" "Q"
Add-QA 16 "Puedo enviarlos a un repositorio remoto" "A"
Add-QA 17 "Ambos apuntan a un commit " "A"
Add-QA 18 "Si hago commit en el branch rama1" "Q"
Add-QA 19 "Solo veré el commit en el historial del branch rama1" "A"
Add-QA 20 "El operador ~ indica" "Q"
Add-QA 21 "Previo sobre mainline" "A"
Add-QA 22 "(sha)^3 funcionaría solo si" "Q"
Add-QA 23 "El commit tiene 3 o más parents" "A"
Add-QA 24 "Un merge resuelto con estrategia fast-forward" "Q"
Add-QA 25 "Implica que solo una de ellas tenía contenido" "A"
Add-QA 26 "Un branch o rama es:" "Q"
Add-QA 27 " Un puntero o variable dirigido a un commit" "A"
Add-QA 28 "Si otro usuario del repositorio crea un branch:" "Q"
Add-QA 29 "Tendré que ejecutar git fetch para visualizar las novedades" "A"

# Match the author's final view/selection state on the sheet.
$ws.Activate()
$ws.Range("A32").Select()
